# MUL vaccine + new table for children vaccines
#
# The "model" choice-list sheet had 60 rows (63-122) holding the per-vaccine
# field names VAC1TIPO..VAC20INF (3 fields x 20 vaccines). That per-vaccine
# table is removed (superseded by a new generic "MUL vaccine" child table),
# which also drops the 60 now-unused shared strings and renumbers every
# other shared-string reference in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

for ($r = 63; $r -le 122; $r++) {
    $ws.Range("A" + $r).Value = ""
    $ws.Range("B" + $r).Value = ""
    $ws.Range("C" + $r).Value = ""
}

# Restore the view: scroll back to the top of the frozen pane and select D4.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D4").Select() | Out-Null
